$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'47.554.65"
$ws.Range("E2").Value = '  +4.80%  '
$ws.Range("D3").Value = "'2.491.20"
$ws.Range("E3").Value = '  +2.76%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").Value = "'322.61"
$ws.Range("E5").Value = '  +1.35%  '
$ws.Range("D6").Value = "'105.31"
$ws.Range("E6").Value = '  +2.05%  '
$ws.Range("D7").Value = "'0.525"
$ws.Range("E7").Value = '  +1.70%  '
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("D9").Value = "'0.544"
$ws.Range("E9").Value = '  +2.89%  '
$ws.Range("D10").Value = "'38.24"
$ws.Range("E10").Value = '  +7.48%  '
$ws.Range("D11").Value = "'0.0813"
$ws.Range("E11").Value = '  +1.10%  '
$ws.Range("E12").Value = '  +1.14%  '
$ws.Range("D13").Value = "'18.32"
$ws.Range("E13").Value = '  +0.79%  '
$ws.Range("D14").Value = "'7.19"
$ws.Range("E14").Value = '  +1.80%  '
$ws.Range("D15").Value = "'2.879.14"
$ws.Range("D16").Value = "'2.483.84"
$ws.Range("E16").Value = '  +1.91%  '
$ws.Range("D17").Value = "'0.848"
$ws.Range("E17").Value = '  +0.59%  '
$ws.Range("D18").Value = "'47.432.51"
$ws.Range("E18").Value = '  +4.70%  '
$ws.Range("E19").Value = '  +4.78%  '
$ws.Range("E20").Value = '  +3.74%  '
$ws.Range("E21").Value = '  +1.55%  '
$ws.Range("D22").Value = "'70.66"
$ws.Range("E22").Value = '  +2.59%  '
$ws.Range("E23").Value = '  +6.21%  '
$ws.Range("D24").Value = "'251.50"
$ws.Range("E24").Value = '  +2.99%  '
$ws.Range("E25").Value = '  +3.48%  '
$ws.Range("D26").Value = "'26.23"
$ws.Range("E26").Value = '  +2.03%  '
$ws.Range("E27").Value = '  +0.01%  '
$ws.Range("E28").Value = '  +4.57%  '
$ws.Range("E29").Value = '  +6.57%  '
$ws.Range("D30").Value = "'35.14"
$ws.Range("E30").Value = '  +6.54%  '
$ws.Range("E31").Value = '  +8.48%  '
$ws.Range("D32").Value = "'49.47"
$ws.Range("E32").Value = '  +0.54%  '
$ws.Range("D33").Value = "'19.67"
$ws.Range("E33").Value = '  -3.03%  '
$ws.Range("E34").Value = '  +3.27%  '
$ws.Range("D35").Value = "'0.0785"
$ws.Range("E35").Value = '  +2.43%  '
$ws.Range("E36").Value = '  +0.13%  '
$ws.Range("E37").Value = '  +5.53%  '
$ws.Range("D38").Value = "'4.63"
$ws.Range("E38").Value = '  +3.84%  '
$ws.Range("D39").Value = "'2.99"
$ws.Range("E39").Value = '  +4.13%  '
$ws.Range("E40").Value = '  +2.11%  '
$ws.Range("E41").Value = '  +1.65%  '
$ws.Range("D42").Value = "'121.63"
$ws.Range("E42").Value = '  -3.61%  '
$ws.Range("D43").Value = "'21.13"
$ws.Range("E43").Value = '  +3.28%  '
$ws.Range("D44").Value = "'0.0298"
$ws.Range("E44").Value = '  +2.51%  '
$ws.Range("D45").Value = "'1.966.44"
$ws.Range("E45").Value = '  +2.06%  '
$ws.Range("D46").Value = "'2.99"
$ws.Range("E46").Value = '  +2.32%  '
$ws.Range("D47").Value = "'2.11"
$ws.Range("E47").Value = '  -0.40%  '
$ws.Range("E48").Value = '  +1.14%  '
$ws.Range("D49").Value = "'9.18"
$ws.Range("E49").Value = '  +0.42%  '
$ws.Range("D50").Value = "'5.27"
$ws.Range("E50").Value = '  +11.74%  '
$ws.Range("D51").Value = "'79.53"
$ws.Range("E51").Value = '  +3.69%  '
